$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 829, shifting row 829 (and everything
# below it, through the old row 870) down by one row -> new row 871.
$ws.Rows.Item(829).Insert()

# Populate the newly inserted row 829 with its data.
# Column A holds a date formatted as text (e.g. "2026/12/29" elsewhere
# in the sheet is stored as a literal string, not a date serial), so
# force the cell to text format before assigning the value to keep
# Excel from auto-converting it into a date number.
$ws.Range("A829").NumberFormat = "@"
$ws.Range("A829").Value = "2026/02/19"
# Drop the number-format override again so the cell ends up with the
# same (default/no explicit) style as its sibling data cells.
$ws.Range("A829").Style = "Normal"

$ws.Range("B829").Value = "木"
$ws.Range("C829").Value = 0
$ws.Range("D829").Value = 201
